# Updates the crypto price/volume snapshot table (rows 2-51) to the latest
# scrape, including the OKB / RenderToken rank swap at rows 43-44.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text so numeric-looking strings (e.g. "166.90", "0.0730") keep their
# exact textual form instead of being parsed into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.079.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.434.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.90"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.02%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.437.76"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.028.94"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.94%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.25"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.01%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.104.49"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.439.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.66"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.48"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.18"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.11"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.19"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.876"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.86"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0730"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.03"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.783.39"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.13%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.89"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.25"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0308"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "333.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.06"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.35"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.53%  "

